$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "addPhoto" translation row (row 37), following the
# existing id / en / th / cn / jp column layout.
$ws.Range("A37").Value = "addPhoto"
$ws.Range("B37").Value = "Add Photo"
$ws.Range("C37").Value = "เพิ่มรูปภาพ"
$ws.Range("D37").Value = "添加照片"
$ws.Range("E37").Value = "写真を追加"

# Update the stored selection to match the author's saved view state.
$ws.Range("H41").Select()
